# Update CRSRCoveredCalls (sheet2) values in column B, rows 11-15,
# and add a new row 16, then make CRSRCoveredCalls the active/selected sheet.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("CRSRShares")
$ws2 = $wb.Worksheets.Item("CRSRCoveredCalls")

# Update existing covered-calls data
$ws2.Range("B11").Value = 11
$ws2.Range("B12").Value = 8
$ws2.Range("B13").Value = 21
$ws2.Range("B14").Value = 22
$ws2.Range("B15").Value = 11

# Append a new row of data
$ws2.Range("A16").Value = 19
$ws2.Range("B16").Value = 0

# Update selection on the CRSRShares sheet (no longer the active tab)
$ws1.Range("C13").Select()

# Activate CRSRCoveredCalls sheet and set its selection
$ws2.Activate()
$ws2.Range("B16").Select()
